# Auto-generated edit script: updates specific leveling-profit cells
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# refreshed Market Board values, per the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 80000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 80000
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = 80000
$ws.Range("N75").Value = -81872
$ws.Range("H78").Value = 80000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 80000
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = 240000
$ws.Range("N78").Value = -249360
$ws.Range("H80").Value = 52490.81
$ws.Range("I80").Value = 1779.9
$ws.Range("J80").Value = 98591.63
$ws.Range("K80").Value = 5339.700000000001
$ws.Range("L80").Value = 295774.89
$ws.Range("M80").Value = -4341.700000000001
$ws.Range("N80").Value = -297770.89
$ws.Range("H83").Value = 52490.81
$ws.Range("I83").Value = 1779.9
$ws.Range("J83").Value = 98591.63
$ws.Range("K83").Value = 16019.1
$ws.Range("L83").Value = 887324.67
$ws.Range("M83").Value = -11027.1
$ws.Range("N83").Value = -897308.67
$ws.Range("H106").Value = 2160.5
$ws.Range("I106").Value = 2160.5
$ws.Range("K106").Value = 2160.5
$ws.Range("M106").Value = -1529.5
$ws.Range("H112").Value = 1035.1
$ws.Range("J112").Value = 983.44446
$ws.Range("L112").Value = 2950.33338
$ws.Range("N112").Value = -5166.33338
$ws.Range("H132").Value = 911.05884
$ws.Range("I132").Value = 929.2
$ws.Range("J132").Value = 775
$ws.Range("K132").Value = 2787.6
$ws.Range("L132").Value = 2325
$ws.Range("M132").Value = -257.6000000000004
$ws.Range("N132").Value = -7385
$ws.Range("H137").Value = 2710.75
$ws.Range("J137").Value = 3331
$ws.Range("L137").Value = 9993
$ws.Range("N137").Value = -15093
$ws.Range("H138").Value = 4026.7532
$ws.Range("J138").Value = 4252.902
$ws.Range("L138").Value = 12758.706
$ws.Range("N138").Value = -23038.706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4263.6
$ws.Range("I61").Value = 1828.6
$ws.Range("J61").Value = 6698.6
$ws.Range("K61").Value = 1828.6
$ws.Range("L61").Value = 6698.6
$ws.Range("M61").Value = -1616.6
$ws.Range("N61").Value = -7122.6
$ws.Range("H102").Value = 1323.6111
$ws.Range("I102").Value = 1323.6111
$ws.Range("K102").Value = 1323.6111
$ws.Range("M102").Value = 298.3888999999999
$ws.Range("H136").Value = 4263.6
$ws.Range("I136").Value = 1828.6
$ws.Range("J136").Value = 6698.6
$ws.Range("K136").Value = 5485.799999999999
$ws.Range("L136").Value = 20095.8
$ws.Range("M136").Value = -2935.799999999999
$ws.Range("N136").Value = -25195.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 700.875
$ws.Range("J94").Value = 599
$ws.Range("L94").Value = 599
$ws.Range("N94").Value = -1501
$ws.Range("H96").Value = 10752.667
$ws.Range("I96").Value = 10752.667
$ws.Range("K96").Value = 10752.667
$ws.Range("M96").Value = -8006.666999999999
$ws.Range("H99").Value = 2936.625
$ws.Range("I99").Value = 2453.5715
$ws.Range("K99").Value = 2453.5715
$ws.Range("M99").Value = -955.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 14897.529
$ws.Range("I22").Value = 179.08333
$ws.Range("K22").Value = 179.08333
$ws.Range("M22").Value = 170.91667
$ws.Range("H31").Value = 7166.5
$ws.Range("I31").Value = 8747.5
$ws.Range("J31").Value = 6903
$ws.Range("K31").Value = 8747.5
$ws.Range("L31").Value = 6903
$ws.Range("M31").Value = -8452.5
$ws.Range("N31").Value = -7493
$ws.Range("H34").Value = 7166.5
$ws.Range("I34").Value = 8747.5
$ws.Range("J34").Value = 6903
$ws.Range("K34").Value = 8747.5
$ws.Range("L34").Value = 6903
$ws.Range("M34").Value = -8545.5
$ws.Range("N34").Value = -7307
$ws.Range("H93").Value = 19900
$ws.Range("I93").Value = 19900
$ws.Range("K93").Value = 19900
$ws.Range("M93").Value = -18028

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 712.2
$ws.Range("I50").Value = 439
$ws.Range("K50").Value = 1317
$ws.Range("M50").Value = -836
$ws.Range("H53").Value = 712.2
$ws.Range("I53").Value = 439
$ws.Range("K53").Value = 1317
$ws.Range("M53").Value = -836
$ws.Range("H97").Value = 5209415.5
$ws.Range("I97").Value = 572.25
$ws.Range("J97").Value = 7813837
$ws.Range("K97").Value = 1716.75
$ws.Range("L97").Value = 23441511
$ws.Range("M97").Value = -1220.75
$ws.Range("N97").Value = -23442503
$ws.Range("H111").Value = 230
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H130").Value = 2367.5
$ws.Range("I130").Value = 1993.3334
$ws.Range("K130").Value = 5980.0002
$ws.Range("M130").Value = -960.0002000000004
$ws.Range("H134").Value = 15034
$ws.Range("J134").Value = 18988
$ws.Range("L134").Value = 56964
$ws.Range("N134").Value = -67104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 113952
$ws.Range("I4").Value = 225000
$ws.Range("J4").Value = 2904
$ws.Range("K4").Value = 225000
$ws.Range("L4").Value = 2904
$ws.Range("M4").Value = -224888
$ws.Range("N4").Value = -3128
$ws.Range("H57").Value = 14000
$ws.Range("I57").Value = 14000
$ws.Range("K57").Value = 14000
$ws.Range("M57").Value = -13180
$ws.Range("H97").Value = 740.3077
$ws.Range("I97").Value = 796.25
$ws.Range("J97").Value = 715.44446
$ws.Range("K97").Value = 796.25
$ws.Range("L97").Value = 715.44446
$ws.Range("M97").Value = -300.25
$ws.Range("N97").Value = -1707.44446
$ws.Range("H132").Value = 3066.3333
$ws.Range("I132").Value = 3049.6667
$ws.Range("K132").Value = 9149.000100000001
$ws.Range("M132").Value = -6619.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1072.9
$ws.Range("J16").Value = 1265.6666
$ws.Range("L16").Value = 1265.6666
$ws.Range("N16").Value = -1605.6666
$ws.Range("H22").Value = 2387.2666
$ws.Range("I22").Value = 2596.625
$ws.Range("J22").Value = 2148
$ws.Range("K22").Value = 2596.625
$ws.Range("L22").Value = 2148
$ws.Range("M22").Value = -2301.625
$ws.Range("N22").Value = -2738
$ws.Range("H27").Value = 2387.2666
$ws.Range("I27").Value = 2596.625
$ws.Range("J27").Value = 2148
$ws.Range("K27").Value = 2596.625
$ws.Range("L27").Value = 2148
$ws.Range("M27").Value = -2489.625
$ws.Range("N27").Value = -2362
$ws.Range("H55").Value = 1060.0952
$ws.Range("I55").Value = 1049.3
$ws.Range("J55").Value = 1069.909
$ws.Range("K55").Value = 1049.3
$ws.Range("L55").Value = 1069.909
$ws.Range("M55").Value = -876.3
$ws.Range("N55").Value = -1415.909
$ws.Range("H122").Value = 4932.108
$ws.Range("I122").Value = 4135.923
$ws.Range("J122").Value = 6814
$ws.Range("K122").Value = 12407.769
$ws.Range("L122").Value = 20442
$ws.Range("M122").Value = -9957.769
$ws.Range("N122").Value = -25342
$ws.Range("H132").Value = 3488.1875
$ws.Range("I132").Value = 2800.6
$ws.Range("K132").Value = 8401.799999999999
$ws.Range("M132").Value = -5871.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9446.75
$ws.Range("I62").Value = 6894.5
$ws.Range("K62").Value = 6894.5
$ws.Range("M62").Value = -6270.5
$ws.Range("H65").Value = 9446.75
$ws.Range("I65").Value = 6894.5
$ws.Range("K65").Value = 34472.5
$ws.Range("M65").Value = -31352.5
$ws.Range("H132").Value = 3949.5
$ws.Range("I132").Value = 3949.5
$ws.Range("K132").Value = 11848.5
$ws.Range("M132").Value = -9318.5
$ws.Range("H136").Value = 2517.4
$ws.Range("I136").Value = 2517.4
$ws.Range("K136").Value = 7552.200000000001
$ws.Range("M136").Value = -5002.200000000001
